# Excel alignment constants
$xlLeft = -4131
$xlTop  = -4160

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Row 17 ("repaymentstrategy") B-cell: swap the value from "Mifos style" to
# the new "Penalties, Fees, Interest, Principal order" option, and restyle it
# to left/top aligned without wrap (a new, distinct cell style).
$cell = $ws.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlTop
$cell.WrapText = $false

# Leave the selection on the edited cell, matching the saved view state.
$ws.Range("B17").Select()
